# Weekly update: insert a new record row for Puerro (Vega Central Mapocho de
# Santiago) above the existing row 34, shifting all subsequent rows down by
# one. Excel's native Rows.Insert() reproduces the row-shift + dimension
# growth seen in the diff (A1:R98 -> A1:R99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 34; rows 34-98 shift down to 35-99.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record.
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44720
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100112005
$ws.Range("G34").Value = "Puerro"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 160
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = 7500
$ws.Range("N34").Value = "`$/paquete 20 unidades"
$ws.Range("O34").Value = "Provincia de Chacabuco"
$ws.Range("P34").Value = 375
$ws.Range("Q34").Value = 20
$ws.Range("R34").Value = "Hortaliza"
